# Apply the "Added fuel2 and efficiency2 to data sheet" edit to both
# worksheets (lowBio, highBio). Each sheet's old rows 4 (EC_OCAES) and
# 5 (EC_VFB) are removed / replaced: the data that used to live in row 6
# (Fuels / MaxActivity / BIO / [PJ]) now lives in row 4, row 5 becomes a
# blank spacer row (formatting only) and row 6 is deleted outright.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: lowBio
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("lowBio")

# Row 4 becomes the "fuel1 / efficiency1" row (previously row 6 data).
$ws1.Range("E4:F4").ClearFormats()
$ws1.Range("A4").Value = "Fuels"
$ws1.Range("B4").Value = "MaxActivity"
$ws1.Range("C4").Value = "BIO"
$ws1.Range("D4").Value = "[PJ]"
$ws1.Range("E4").Value = 52.6
$ws1.Range("F4").Value = "constant"
$ws1.Range("G4").Value = 52.6
$ws1.Range("H4").Value = 105.2
$ws1.Range("J4").ClearContents()

# Row 5 becomes a blank spacer row (keep cell formatting, drop content).
$ws1.Range("A5:J5").ClearContents()

# Old row 6 is removed entirely.
$ws1.Rows("6").Delete()

# ---------------------------------------------------------------------
# Sheet 2: highBio
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("highBio")

# EC_DAC high value updated.
$ws2.Range("E3").Value = 3704

# Row 4 becomes the "fuel2 / efficiency2" row (previously row 6 data,
# now a plain value instead of a formula).
$ws2.Range("E4:F4").ClearFormats()
$ws2.Range("A4").Value = "Fuels"
$ws2.Range("B4").Value = "MaxActivity"
$ws2.Range("C4").Value = "BIO"
$ws2.Range("D4").Value = "[PJ]"
$ws2.Range("E4").Value = 105.2
$ws2.Range("F4").Value = "constant"
$ws2.Range("G4").Value = 52.6
$ws2.Range("H4").Value = 105.2
$ws2.Range("I4:J4").Clear()

# Row 5 becomes a blank spacer row (keep cell formatting, drop content).
$ws2.Range("A5:J5").ClearContents()

# Old row 6 is removed entirely.
$ws2.Rows("6").Delete()

# ---------------------------------------------------------------------
# Restore the selections shown in the edited file (highBio selected
# first so lowBio ends up the active/tab-selected sheet, as before).
# ---------------------------------------------------------------------
$null = $ws2.Range("E5").Select()
$null = $ws1.Range("F10").Select()
